# Auto-generated edit script: updates cryptos worksheet values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.712.03'
$ws.Range("E2").Value = '  +2.46%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.368.84'
$ws.Range("E3").Value = '  +6.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.11'
$ws.Range("E5").Value = '  +6.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.58'
$ws.Range("E6").Value = '  -1.88%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.643'
$ws.Range("E7").Value = '  +2.69%  '

# Row 8
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.644'
$ws.Range("E9").Value = '  +5.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.34'
$ws.Range("E10").Value = '  -2.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0942'
$ws.Range("E11").Value = '  +2.66%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.82'
$ws.Range("E12").Value = '  -0.69%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.05'
$ws.Range("E13").Value = '  +3.52%  '

# Row 14
$ws.Range("E14").Value = '  +2.09%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.54'
$ws.Range("E15").Value = '  +9.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.726.95'
$ws.Range("E16").Value = '  +6.56%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.367.85'
$ws.Range("E17").Value = '  +6.47%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.702.68'
$ws.Range("E18").Value = '  +2.94%  '

# Row 19
$ws.Range("E19").Value = '  +2.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.29'
$ws.Range("E20").Value = '  -1.80%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.37'
$ws.Range("E21").Value = '  +3.74%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.46'
$ws.Range("E22").Value = '  -1.64%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.56'
$ws.Range("E23").Value = '  +7.98%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '257.66'
$ws.Range("E24").Value = '  +12.27%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.33'
$ws.Range("E25").Value = '  +0.73%  '

# Row 26
$ws.Range("E26").Value = '  +3.88%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.40'
$ws.Range("E28").Value = '  +2.81%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +0.93%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.69'
$ws.Range("E30").Value = '  +7.47%  '

# Row 31
$ws.Range("B31").Value = 'WEMIXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.20'
$ws.Range("E31").Value = '  -1.21%  '

# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.62'
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0929'
$ws.Range("E33").Value = '  +3.73%  '

# Row 34
$ws.Range("E34").Value = '  +5.61%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  +5.21%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.99'
$ws.Range("E36").Value = '  -4.51%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("E37").Value = '  -4.00%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0377'
$ws.Range("E38").Value = '  -0.66%  '

# Row 39
$ws.Range("E39").Value = '  -0.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.84'
$ws.Range("E40").Value = '  +16.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.51'
$ws.Range("E41").Value = '  +13.77%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.11'
$ws.Range("E42").Value = '  -0.32%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.234'
$ws.Range("E43").Value = '  -1.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.88'
$ws.Range("E44").Value = '  +0.66%  '

# Row 45
$ws.Range("E45").Value = '  +0.21%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.67'
$ws.Range("E46").Value = '  +3.49%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.10'
$ws.Range("E47").Value = '  +8.37%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.35'
$ws.Range("E48").Value = '  +10.11%  '

# Row 49
$ws.Range("E49").Value = '  -0.50%  '

# Row 50
$ws.Range("E50").Value = '  +2.72%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.480'
$ws.Range("E51").Value = '  +8.59%  '

Write-Output "Updated cryptos list values"
